$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1160.5
$ws.Range("I43").Value = 802.1667
$ws.Range("J43").Value = 1429.25
$ws.Range("K43").Value = 802.1667
$ws.Range("L43").Value = 1429.25
$ws.Range("M43").Value = -733.1667
$ws.Range("N43").Value = -1567.25

$ws.Range("H49").Value = 3500
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 4125
$ws.Range("K49").Value = 3000
$ws.Range("L49").Value = 12375
$ws.Range("M49").Value = -2864
$ws.Range("N49").Value = -12647

$ws.Range("H74").Value = 3998.853
$ws.Range("I74").Value = 4342.1763
$ws.Range("J74").Value = 3655.5293
$ws.Range("K74").Value = 4342.1763
$ws.Range("L74").Value = 3655.5293
$ws.Range("M74").Value = -3406.1763
$ws.Range("N74").Value = -5527.5293

$ws.Range("H76").Value = 79670.8
$ws.Range("I76").Value = 126172.89
$ws.Range("J76").Value = 3576.4546
$ws.Range("K76").Value = 126172.89
$ws.Range("L76").Value = 3576.4546
$ws.Range("M76").Value = -125857.89
$ws.Range("N76").Value = -4206.4546

$ws.Range("H77").Value = 3998.853
$ws.Range("I77").Value = 4342.1763
$ws.Range("J77").Value = 3655.5293
$ws.Range("K77").Value = 21710.8815
$ws.Range("L77").Value = 18277.6465
$ws.Range("M77").Value = -17030.8815
$ws.Range("N77").Value = -27637.6465

$ws.Range("H79").Value = 79670.8
$ws.Range("I79").Value = 126172.89
$ws.Range("J79").Value = 3576.4546
$ws.Range("K79").Value = 126172.89
$ws.Range("L79").Value = 3576.4546
$ws.Range("M79").Value = -125080.89
$ws.Range("N79").Value = -5760.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1018
$ws.Range("I2").Value = 925.2381
$ws.Range("J2").Value = 1147.8667
$ws.Range("K2").Value = 925.2381
$ws.Range("L2").Value = 1147.8667
$ws.Range("M2").Value = -812.2381
$ws.Range("N2").Value = -1373.8667

$ws.Range("H75").Value = 38000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 38000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 38000
$ws.Range("N75").Value = -39748

$ws.Range("H78").Value = 38000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 38000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 114000
$ws.Range("N78").Value = -122736

$ws.Range("H88").Value = 125101380
$ws.Range("I88").Value = 1778
$ws.Range("J88").Value = 250200980
$ws.Range("K88").Value = 1778
$ws.Range("L88").Value = 250200980
$ws.Range("M88").Value = -1372
$ws.Range("N88").Value = -250201792

$ws.Range("H91").Value = 125101380
$ws.Range("I91").Value = 1778
$ws.Range("J91").Value = 250200980
$ws.Range("K91").Value = 1778
$ws.Range("L91").Value = 250200980
$ws.Range("M91").Value = -374
$ws.Range("N91").Value = -250203788

$ws.Range("H97").Value = 2291.9048
$ws.Range("I97").Value = 2512.7778
$ws.Range("J97").Value = 966.6667
$ws.Range("K97").Value = 2512.7778
$ws.Range("L97").Value = 966.6667
$ws.Range("M97").Value = -2016.7778
$ws.Range("N97").Value = -1958.6667

$ws.Range("H116").Value = 1018
$ws.Range("I116").Value = 925.2381
$ws.Range("J116").Value = 1147.8667
$ws.Range("K116").Value = 925.2381
$ws.Range("L116").Value = 1147.8667
$ws.Range("M116").Value = 1368.7619
$ws.Range("N116").Value = -5735.8667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1018
$ws.Range("I3").Value = 925.2381
$ws.Range("J3").Value = 1147.8667
$ws.Range("K3").Value = 925.2381
$ws.Range("L3").Value = 1147.8667
$ws.Range("M3").Value = -811.2381
$ws.Range("N3").Value = -1375.8667

$ws.Range("H86").Value = 7144771.5
$ws.Range("I86").Value = 9092750
$ws.Range("J86").Value = 2183.1667
$ws.Range("K86").Value = 9092750
$ws.Range("L86").Value = 2183.1667
$ws.Range("M86").Value = -9091627
$ws.Range("N86").Value = -4429.1667

$ws.Range("H89").Value = 7144771.5
$ws.Range("I89").Value = 9092750
$ws.Range("J89").Value = 2183.1667
$ws.Range("K89").Value = 45463750
$ws.Range("L89").Value = 10915.8335
$ws.Range("M89").Value = -45458134
$ws.Range("N89").Value = -22147.8335

$ws.Range("H94").Value = 9076.84
$ws.Range("I94").Value = 961.8125
$ws.Range("J94").Value = 23503.555
$ws.Range("K94").Value = 961.8125
$ws.Range("L94").Value = 23503.555
$ws.Range("M94").Value = -510.8125
$ws.Range("N94").Value = -24405.555

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 20039.666
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 20039.666
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 20039.666
$ws.Range("N59").Value = -22329.666

$ws.Range("H62").Value = 47622028
$ws.Range("I62").Value = 2818.5
$ws.Range("J62").Value = 90912216
$ws.Range("K62").Value = 2818.5
$ws.Range("L62").Value = 90912216
$ws.Range("M62").Value = -2194.5
$ws.Range("N62").Value = -90913464

$ws.Range("H65").Value = 47622028
$ws.Range("I65").Value = 2818.5
$ws.Range("J65").Value = 90912216
$ws.Range("K65").Value = 14092.5
$ws.Range("L65").Value = 454561080
$ws.Range("M65").Value = -10972.5
$ws.Range("N65").Value = -454567320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H63").Value = 8571.817999999999
$ws.Range("I63").Value = 1430
$ws.Range("J63").Value = 11250
$ws.Range("K63").Value = 4290
$ws.Range("L63").Value = 33750
$ws.Range("M63").Value = -3541
$ws.Range("N63").Value = -35248

$ws.Range("H66").Value = 8571.817999999999
$ws.Range("I66").Value = 1430
$ws.Range("J66").Value = 11250
$ws.Range("K66").Value = 12870
$ws.Range("L66").Value = 101250
$ws.Range("M66").Value = -9126
$ws.Range("N66").Value = -108738

$ws.Range("H76").Value = 6070.7144
$ws.Range("I76").Value = 1663.3334
$ws.Range("J76").Value = 7272.727
$ws.Range("K76").Value = 4990.0002
$ws.Range("L76").Value = 21818.181
$ws.Range("M76").Value = -4607.0002
$ws.Range("N76").Value = -22584.181

$ws.Range("H79").Value = 6070.7144
$ws.Range("I79").Value = 1663.3334
$ws.Range("J79").Value = 7272.727
$ws.Range("K79").Value = 4990.0002
$ws.Range("L79").Value = 21818.181
$ws.Range("M79").Value = -3664.0002
$ws.Range("N79").Value = -24470.181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3817.8125
$ws.Range("I80").Value = 6281
$ws.Range("J80").Value = 2698.182
$ws.Range("K80").Value = 6281
$ws.Range("L80").Value = 2698.182
$ws.Range("M80").Value = -5283
$ws.Range("N80").Value = -4694.182

$ws.Range("H83").Value = 3817.8125
$ws.Range("I83").Value = 6281
$ws.Range("J83").Value = 2698.182
$ws.Range("K83").Value = 31405
$ws.Range("L83").Value = 13490.91
$ws.Range("M83").Value = -26413
$ws.Range("N83").Value = -23474.91

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1581.3889
$ws.Range("I7").Value = 1249.2354
$ws.Range("J7").Value = 1878.579
$ws.Range("K7").Value = 1249.2354
$ws.Range("L7").Value = 1878.579
$ws.Range("M7").Value = -1137.2354
$ws.Range("N7").Value = -2102.579

$ws.Range("H82").Value = 6062188.5
$ws.Range("I82").Value = 9091803
$ws.Range("J82").Value = 2960
$ws.Range("K82").Value = 9091803
$ws.Range("L82").Value = 2960
$ws.Range("M82").Value = -9091442
$ws.Range("N82").Value = -3682

$ws.Range("H85").Value = 6062188.5
$ws.Range("I85").Value = 9091803
$ws.Range("J85").Value = 2960
$ws.Range("K85").Value = 9091803
$ws.Range("L85").Value = 2960
$ws.Range("M85").Value = -9090555
$ws.Range("N85").Value = -5456

$ws.Range("H126").Value = 1581.3889
$ws.Range("I126").Value = 1249.2354
$ws.Range("J126").Value = 1878.579
$ws.Range("K126").Value = 3747.7062
$ws.Range("L126").Value = 5635.737
$ws.Range("M126").Value = -1277.7062
$ws.Range("N126").Value = -10575.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1158.9491
$ws.Range("I132").Value = 624.1064
$ws.Range("J132").Value = 3253.75
$ws.Range("K132").Value = 1872.3192
$ws.Range("L132").Value = 9761.25
$ws.Range("M132").Value = 657.6808000000001
$ws.Range("N132").Value = -14821.25
